# Fix player loop bug: update TOTAL RUNS (col B) and WICKETS (col C) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 45
    "B3" = 167
    "B4" = 103
    "B5" = 90
    "B6" = 76
    "B7" = 110
    "B8" = 14
    "C8" = 8
    "B9" = 85
    "C9" = 9
    "B10" = 70
    "C10" = 13
    "B11" = 37
    "C11" = 11
    "B12" = 39
    "C12" = 7
    "B24" = 72
    "B25" = 57
    "B26" = 49
    "B27" = 198
    "B28" = 73
    "B29" = 86
    "B30" = 122
    "C30" = 4
    "B31" = 115
    "C31" = 7
    "B32" = 71
    "C32" = 6
    "B33" = 65
    "C33" = 12
    "B34" = 44
    "C34" = 14
    "B35" = 26
    "B36" = 116
    "B37" = 91
    "B39" = 33
    "B40" = 136
    "B41" = 46
    "C41" = 4
    "B42" = 27
    "C42" = 6
    "C43" = 11
    "B44" = 85
    "C44" = 10
    "B45" = 19
    "C45" = 9
    "B46" = 83
    "B47" = 65
    "B48" = 91
    "B49" = 25
    "B50" = 78
    "B51" = 126
    "B52" = 197
    "C52" = 10
    "B53" = 76
    "C53" = 8
    "B54" = 54
    "C54" = 10
    "B55" = 24
    "C55" = 13
    "B56" = 14
    "C56" = 9
    "B57" = 179
    "B58" = 136
    "B59" = 142
    "B60" = 237
    "B61" = 109
    "B62" = 115
    "B63" = 58
    "C63" = 16
    "B64" = 41
    "C64" = 13
    "B65" = 126
    "C65" = 12
    "B66" = 65
    "C66" = 17
    "B67" = 121
    "C67" = 17
    "B68" = 44
    "B69" = 72
    "B70" = 55
    "B71" = 65
    "B72" = 26
    "B74" = 134
    "C74" = 7
    "B75" = 37
    "C75" = 8
    "B76" = 87
    "C76" = 9
    "B77" = 60
    "C77" = 5
    "B79" = 184
    "B80" = 58
    "B81" = 26
    "B82" = 104
    "B83" = 75
    "B85" = 49
    "C85" = 2
    "B86" = 56
    "C86" = 6
    "C87" = 9
    "B88" = 34
    "C88" = 14
    "B89" = 27
    "C89" = 15
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
